$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1 ("Naam" row value): keep the bold font, but also give it the same
# highlight fill that the rest of column B already uses. Copy the fill from
# a cell that already has it (B2), then restore the bold font that the
# format-paste just clobbered.
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Font.Bold = $true

# --- B6 (Tabel 24 notes): it was using a different ("theme 7") highlight
# color than every other note cell; switch it to the standard ("theme 9")
# highlight used everywhere else, while keeping its wrap-text alignment.
# Copy that combination of formatting from B3, which already has it.
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- B6 text: drop the leftover "TODO Mat54" notes on the 3524 / 3526 lines
$ws.Range("B6").Value = "Tabel 24a*:`n[x] 5221`n[x] 9621`n[x] 923 >?`nTabel 24b:`n[x] 5224`n[x] 3524`n[n] 9624`n[x] 924 <?`n[x] 9626`n[x] 5226`n[x] 3526"

# --- View state: the sheet had scrolled so row 6 was at the top with B19
# selected; now it is scrolled back to the top with B4 selected instead.
[void]$ws.Range("B4").Select()
